{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that preceded it) that followed the\n// bibliography entry ending in \"Janeiro: Editora Interci\u00eancia , 2004.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography paragraph that anchors the block we need to drop.\nconst anchorText = \"Janeiro: Editora Interci\u00eancia , 2004.\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const blankText = \"\";\n  const jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n  const copyrightPrefix = \"\u00a9 2020\";\n\n  // The three paragraphs immediately following the anchor are: a blank\n  // paragraph, the \"Ver no Jupiter...\" line, and the \"\u00a9 2020...\" line.\n  // Verify and delete them (in reverse order so indices stay valid).\n  const toDelete = [];\n  if (\n    anchorIndex + 3 < items.length &&\n    items[anchorIndex + 1].text.trim() === blankText &&\n    items[anchorIndex + 2].text.trim() === jupiterText &&\n    items[anchorIndex + 3].text.trim().startsWith(copyrightPrefix)\n  ) {\n    toDelete.push(items[anchorIndex + 1], items[anchorIndex + 2], items[anchorIndex + 3]);\n  }\n\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that preceded it) that followed the\n# bibliography entry ending in \"Janeiro: Editora Interci\u00eancia , 2004.\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Janeiro: Editora Interci\u00eancia , 2004.\"\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($anchorText)\n\nif ($found) {\n    # Resolve the 1-based Paragraphs() index of the paragraph containing the\n    # found text by counting paragraph marks before it.\n    $precedingRange = $d.Range(0, $searchRange.Start)\n    $anchorIndex = $precedingRange.Paragraphs.Count + 1\n\n    if (($anchorIndex + 3) -le $d.Paragraphs.Count) {\n        $blankPara = $d.Paragraphs.Item($anchorIndex + 1)\n        $jupiterPara = $d.Paragraphs.Item($anchorIndex + 2)\n        $copyrightPara = $d.Paragraphs.Item($anchorIndex + 3)\n\n        $blankText = $blankPara.Range.Text.Trim()\n        $jupiterText = $jupiterPara.Range.Text.Trim()\n        $copyrightText = $copyrightPara.Range.Text.Trim()\n\n        if ($blankText -eq \"\" -and $jupiterText -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\" -and $copyrightText.StartsWith(\"\u00a9 2020\")) {\n            $deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)\n            $deleteRange.Delete()\n        }\n    }\n}\n"}
